$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '63.875.12'
Set-TextValue 'E2' '  +0.28%  '
Set-TextValue 'D3' '2.752.90'
Set-TextValue 'E3' '  +0.69%  '
Set-TextValue 'E4' '  -0.19%  '
Set-TextValue 'D5' '579.63'
Set-TextValue 'E5' '  -2.28%  '
Set-TextValue 'D6' '159.52'
Set-TextValue 'E6' '  +4.32%  '
Set-TextValue 'E7' '  +0.39%  '
Set-TextValue 'E8' '  +0.30%  '
Set-TextValue 'E9' '  -0.80%  '
Set-TextValue 'D10' '0.393'
Set-TextValue 'E10' '  +1.03%  '
Set-TextValue 'D11' '5.71'
Set-TextValue 'E11' '  -15.13%  '
Set-TextValue 'E12' '  +0.30%  '
Set-TextValue 'D13' '3.238.92'
Set-TextValue 'E13' '  +0.31%  '
Set-TextValue 'D14' '27.00'
Set-TextValue 'E14' '  +1.52%  '
Set-TextValue 'D15' '63.824.68'
Set-TextValue 'E15' '  +0.32%  '
Set-TextValue 'E16' '  +1.67%  '
Set-TextValue 'D17' '2.754.25'
Set-TextValue 'E17' '  -0.28%  '
Set-TextValue 'D18' '12.28'
Set-TextValue 'E18' '  +1.89%  '
Set-TextValue 'E19' '  +1.07%  '
Set-TextValue 'D20' '362.27'
Set-TextValue 'E20' '  -0.93%  '
Set-TextValue 'E21' '  -1.78%  '
Set-TextValue 'E22' '  +6.00%  '
Set-TextValue 'E23' '  +0.17%  '
Set-TextValue 'D24' '66.39'
Set-TextValue 'E24' '  +0.68%  '
Set-TextValue 'E25' '  +2.66%  '
Set-TextValue 'D26' '8.67'
Set-TextValue 'E26' '  +0.03%  '
Set-TextValue 'E27' '  +0.45%  '
Set-TextValue 'D28' '0.0₃0942'
Set-TextValue 'E28' '  +3.51%  '
Set-TextValue 'E29' '  -1.25%  '
Set-TextValue 'E30' '  +0.05%  '
Set-TextValue 'E31' '  +4.73%  '
Set-TextValue 'D32' '168.71'
Set-TextValue 'D34' '20.59'
Set-TextValue 'E34' '  -0.22%  '
Set-TextValue 'E35' '  +3.63%  '
Set-TextValue 'E36' '  +2.74%  '
Set-TextValue 'E37' '  +1.78%  '
Set-TextValue 'E38' '  +0.32%  '
Set-TextValue 'E39' '  -0.30%  '
Set-TextValue 'D40' '6.17'
Set-TextValue 'E40' '  +10.18%  '
Set-TextValue 'D41' '332.29'
Set-TextValue 'E41' '  -4.56%  '
Set-TextValue 'D42' '39.56'
Set-TextValue 'E42' '  +1.57%  '
Set-TextValue 'D43' '22.09'
Set-TextValue 'E43' '  -0.24%  '
Set-TextValue 'B44' 'InjectiveProtocol'
Set-TextValue 'C44' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D44' '22.04'
Set-TextValue 'E44' '  -0.86%  '
Set-TextValue 'B45' 'Hedera'
Set-TextValue 'C45' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D45' '0.0600'
Set-TextValue 'E45' '  +1.22%  '
Set-TextValue 'B46' 'Mantle'
Set-TextValue 'C46' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D46' '0.641'
Set-TextValue 'E46' '  -0.90%  '
Set-TextValue 'B47' 'VeChain'
Set-TextValue 'C47' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D47' '0.0258'
Set-TextValue 'E47' '  -0.21%  '
Set-TextValue 'D48' '136.98'
Set-TextValue 'E48' '  -4.45%  '
Set-TextValue 'E49' '  +0.70%  '
Set-TextValue 'E50' '  +0.52%  '
Set-TextValue 'D51' '11.07'
Set-TextValue 'E51' '  +1.02%  '
